$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix a data-entry bug caught by the test suite ---
# A6 (DNI 27775770) had been entered as text; it must be a real number
# like every other DNI in column A.
$ws.Range("A6").Value = 27775770

# --- New check-in record: Majolli, Facundo (row 7) ---
$ws.Range("A7").Value = 44189151
$ws.Range("B7").Value = "Majolli"
$ws.Range("C7").Value = "Facundo"
# Fecha/Hora are stored as plain text, not real dates/times, so force text
# with a leading apostrophe to stop Excel from auto-converting "2023-11-13"
# into a date serial number; then reset the style so no extra number
# format sticks to the cell.
$ws.Range("D7").Value = "'2023-11-13"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "08:26:17"

# --- New check-in record: Reynoso, Anahi (row 8) ---
# DNI on this row is stored as text (matches the source diff) rather than
# a number, so also force it with a leading apostrophe.
$ws.Range("A8").Value = "'31949304"
$ws.Range("A8").Style = "Normal"
$ws.Range("B8").Value = "Reynoso"
$ws.Range("C8").Value = "Anahi"
$ws.Range("D8").Value = "'2023-11-13"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "08:34:09"
